# Uren Registratie Game-Lab-2.1 - update weekly hour entries
# (AIEnemy Abstract af, Asset list en uren registratie)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 5 block (row 5 = Woensdag): 2 -> 4 hours for everybody
$ws.Range("B5:G5").Value = 4

# Week 6 block (row 13 = Woensdag): 2 -> 4, but Dinsdag(D13) stays 0 (absent)
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 4

# Week 7 block (row 23 = Vrijdag): 4 -> 8
$ws.Range("B23:G23").Value = 8

# Week 8 block (row 31 = Vrijdag): 4 -> 8
$ws.Range("B31:G31").Value = 8

# Week 9 block (row 42 = Maandag): 4 -> 8
$ws.Range("B42:G42").Value = 8

# Week 10 block (row 52 = Woensdag): 4 -> 6
$ws.Range("B52:G52").Value = 6

# Week 11 block (row 58 = Maandag): 4 -> 8
$ws.Range("B58:G58").Value = 8

# Week 12 block (row 70 = Vrijdag): 4 -> 8, but C70 (Dinsdag) was 0 -> 4
$ws.Range("B70").Value = 8
$ws.Range("C70").Value = 4
$ws.Range("D70:G70").Value = 8

# Week 13 block (row 75 = Dinsdag): B75 stays 4; C75 4->8, D75 0->4, E75:G75 4->8
$ws.Range("C75").Value = 8
$ws.Range("D75").Value = 4
$ws.Range("E75:G75").Value = 8

# Week 14 block (row 84 = Woensdag): B84 stays 2; C84 2->6, D84 0->4, E84:G84 2->6
$ws.Range("C84").Value = 6
$ws.Range("D84").Value = 4
$ws.Range("E84:G84").Value = 6

# Week 15 block (row 92 = Woensdag): previously empty, now filled in
$ws.Range("B92").Value = 6
$ws.Range("C92").Value = 6
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 6
$ws.Range("F92").Value = 6
$ws.Range("G92").Value = 6

# Mirror the workbook's convention: cells with logged hours get a blue
# fill, a day with 0 hours (absence) gets a red fill (B-column never
# gets shaded in this sheet, matching the rest of the table).
$ws.Range("C92").Interior.Color = 12611584
$ws.Range("E92").Interior.Color = 12611584
$ws.Range("F92").Interior.Color = 12611584
$ws.Range("G92").Interior.Color = 12611584
$ws.Range("D92").Interior.Color = 255

# Update the view state: scrolled position and active selection
$ws.Range("J92").Select()
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
